# Update published Fonds de solidarite data (2022-06-07 refresh)
# For each listed row, update column C (nombre_aides) and column E (montant_total)
# to their new values. Column D (nombre_entreprises) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 10;  C = 278200;  E = 1752278885 },
    @{ Row = 39;  C = 31689;   E = 63371159 },
    @{ Row = 56;  C = 33210;   E = 108834427 },
    @{ Row = 85;  C = 10749;   E = 47053687 },
    @{ Row = 134; C = 5680;    E = 17175156 },
    @{ Row = 143; C = 2070;    E = 3754697 },
    @{ Row = 167; C = 101535;  E = 195051292 },
    @{ Row = 169; C = 562618;  E = 1285077512 },
    @{ Row = 170; C = 367437;  E = 2846599233 },
    @{ Row = 171; C = 115173;  E = 447158183 },
    @{ Row = 174; C = 357262;  E = 1018669995 },
    @{ Row = 175; C = 125563;  E = 813587329 },
    @{ Row = 177; C = 96764;   E = 174782722 },
    @{ Row = 179; C = 235730;  E = 812758054 },
    @{ Row = 186; C = 21936;   E = 40105992 },
    @{ Row = 255; C = 141367;  E = 414361455 },
    @{ Row = 267; C = 84978;   E = 156522342 },
    @{ Row = 311; C = 190856;  E = 586445221 },
    @{ Row = 323; C = 94724;   E = 178792585 }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("E" + $u.Row).Value = $u.E
}
